$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 58

$ws.Cells.Item($row, 1).Value = "EVDFZI"
$ws.Cells.Item($row, 2).Value = "Fusible de chip de reinicio de tambor"
$ws.Cells.Item($row, 3).Value = "3300 3400 3600 5600 5700 5650 5750 5500 5800 5900 5850 5950 8600 8800 9600 9800 9650 C9850 4010D 431D B411 431 432 332 363 MB441 MB451 461 491"
$ws.Cells.Item($row, 4).Value = 10000
$ws.Cells.Item($row, 5).Value = 50000
$ws.Cells.Item($row, 6).Value = 19
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E58-D58)*G58"
$ws.Cells.Item($row, 9).Formula = "=D58*F58"
$ws.Cells.Item($row, 10).Value = 190000
